# Jun's file updates for all IO data and others
#
# Updates the GDP Growth Rates control-settings workbook:
#  - Data sheet: refresh the Real GDP figures (July STEO -> September STEO)
#    for 2020/2021 and let the dependent formulas recalc automatically.
#  - About sheet: refresh the narrative text referencing the STEO vintage
#    and the "as of" date for the pandemic GDP impact data.
#  - Refresh the saved cell selections on the About and Data sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("Data")

# --- Text updates -----------------------------------------------------
# Set the Data sheet's STEO label first so new shared strings are appended
# in the same order as the target workbook ("September STEO" before the
# About-sheet strings).
$ws2.Range("A3").Value = "September STEO"

$ws1.Range("B6").Value = "January 2020 and September 2020"
$ws1.Range("A28").Value = "SARS-CoV-2 pandemic.  It uses the latest data available as of September 9,"

# --- Data updates -------------------------------------------------------
# Updated Real GDP figures (2020 and 2021) from the September STEO.
# Dependent formulas (Data!C8, Data!D8, Data!B12, GDPGR-alternate!B2) will
# recalculate automatically.
$ws2.Range("C3").Value = 18168
$ws2.Range("D3").Value = 18726

# --- Selection / view updates -------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("D4").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A29").Select() | Out-Null
